$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Full set of values for the 12 rows that end up in the sheet.
# row -> (colA text, colB number, colC text-or-$null)
$json1 = '{"type":"https://tools.ietf.org/html/rfc7231#section-6.5.13","title":"Unsupported Media Type","status":415,"traceId":"|ed7b4edd-478dab1ee0599a51."}'
$json2 = '{"type":"https://tools.ietf.org/html/rfc7231#section-6.5.13","title":"Unsupported Media Type","status":415,"traceId":"|ed7b4ede-478dab1ee0599a51."}'
$json3 = '{"type":"https://tools.ietf.org/html/rfc7231#section-6.5.13","title":"Unsupported Media Type","status":415,"traceId":"|ed7b4edf-478dab1ee0599a51."}'
$json4 = '{"type":"https://tools.ietf.org/html/rfc7231#section-6.5.13","title":"Unsupported Media Type","status":415,"traceId":"|ed7b4ee3-478dab1ee0599a51."}'
$json5 = '{"type":"https://tools.ietf.org/html/rfc7231#section-6.5.13","title":"Unsupported Media Type","status":415,"traceId":"|ed7b4ee4-478dab1ee0599a51."}'

$rows = @{
    1  = @("text", 1, $null)
    2  = @("GEThttps://localhost:44393/api/notifications/{id}", 204, $null)
    3  = @("GEThttps://localhost:44393/api/notifications", 200, "[]")
    4  = @("POSThttps://localhost:44393/api/notifications/new", 415, $json1)
    5  = @("PATCHhttps://localhost:44393/api/notifications/seen/{id}", 415, $json1)
    6  = @("POSThttps://localhost:44393/api/account/register", 415, $json2)
    7  = @("PUThttps://localhost:44393/api/account/login", 415, $json3)
    8  = @("GEThttps://localhost:44393/api/account/Token/Valid", 200, '"Yep, still valid."')
    9  = @("GEThttps://localhost:44393/api/vices/mine", 200, "[]")
    10 = @("GEThttps://localhost:44393/api/vices", 200, '[{"name":"Bautura","viceId":"1"},{"name":"Mancare","viceId":"2"},{"name":"Tigari","viceId":"3"}]')
    11 = @("DELETEhttps://localhost:44393/api/vices", 415, $json4)
    12 = @("PUThttps://localhost:44393/api/vices/updateVices", 415, $json5)
}

# Write cells in the exact order the originating automation used, so that the
# resulting shared-strings table lines up value-for-value with the source
# workbook (Excel/iron_native assigns shared-string ids in first-use order).
# Row 2 (A,B) and row 3 (A) were populated first, then a header/summary row
# was inserted at row 1 (A,B), after which the remaining cells were filled in
# row-by-row from row 3's C cell through row 12.
$ws.Cells.Item(2, 1).Value = $rows[2][0]
$ws.Cells.Item(2, 2).Value = $rows[2][1]

$ws.Cells.Item(3, 1).Value = $rows[3][0]

$ws.Cells.Item(1, 1).Value = $rows[1][0]
$ws.Cells.Item(1, 2).Value = $rows[1][1]

$ws.Cells.Item(3, 2).Value = $rows[3][1]
$ws.Cells.Item(3, 3).Value = $rows[3][2]

for ($r = 4; $r -le 12; $r++) {
    $entry = $rows[$r]
    $ws.Cells.Item($r, 1).Value = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    if ($entry[2] -ne $null) {
        $ws.Cells.Item($r, 3).Value = $entry[2]
    }
}
